# P202 Add Column ( MARKET_FORECAST )
#
# Adds a new "Market Forecast" column (AH) to the P202-By-Season template:
#   - AH1 header  = "Market Forecast"
#   - AH2 value   = "&=result.MARKET_FORECAST"
#   - formats the new column with the 0_);[Red](0) number format
#   - widens column AH and extends the AutoFilter + _FilterDatabase
#     defined name so the new column is included
#   - restores the recorded selection state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column content --------------------------------------------------
$ws.Range("AH1").Value = "Market Forecast"
$ws.Range("AH2").Value = "&=result.MARKET_FORECAST"

# Number format used by the column (matches the pattern already used for
# the other numeric-looking result columns on this sheet).
$ws.Range("AH1").NumberFormat = "0_);[Red]\(0\)"
$ws.Range("AH2").NumberFormat = "0_);[Red]\(0\)"

# Column width for the new column.
$ws.Columns.Item(34).ColumnWidth = 14.65

# --- Extend the AutoFilter to cover the new column ------------------------
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:AH1").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$wb.Names.Item("Sheet 1!_FilterDatabase").RefersTo = "='Sheet 1'!`$A`$1:`$AH`$1"

# --- Restore the recorded selection/view state -----------------------------
$null = $ws.Range("AD5").Select()
